# Update the "Förändrad" (Changed) date column (C) for all data rows
# on the active worksheet from serial date 45190 to serial date 45192.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 302
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45190) {
        $cell.Value2 = 45192
    }
}
